# "clean data to 321"
# Extends the subjects data table (Sheet1) with 7 new subjects (315-321),
# each contributing the usual 4 rows (minWord, minPseudo, symp, gob) in
# column order, mirroring the existing block structure of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout (column order) for every subject block: minWord, minPseudo, symp, gob.
$rowTestCodes = @("minWord", "minPseudo", "symp", "gob")

# New subject numbers with their per-test_code values.
$newSubjects = @(
    @{ Subj = 315; minWord = "v7hpnfwykc"; minPseudo = "zvtxod636o"; symp = "ekge7117c0"; gob = "mbse2ual93" },
    @{ Subj = 316; minWord = "ellh3f7r6f"; minPseudo = "dyhqkv4y90"; symp = "b0n958l42w"; gob = "vfdk659lm2" },
    @{ Subj = 317; minWord = "g3h6mfmlh9"; minPseudo = "8jm4ktdp2v"; symp = "bvzbcex95y"; gob = "dmcfba923s" },
    @{ Subj = 318; minWord = "a8fzyae532"; minPseudo = "6ee8c6wkk7"; symp = "7xjxwvk39u"; gob = "sedpugnmk5" },
    @{ Subj = 319; minWord = "z4rx8bqocr"; minPseudo = "qgecou8jlo"; symp = "ma2hmr7lqa"; gob = "37uuqlft2a" },
    @{ Subj = 320; minWord = "50ppjg6by9"; minPseudo = "b2jojnaqxl"; symp = "o8rq1gza2t"; gob = "6tavb9lkp3" },
    @{ Subj = 321; minWord = "s5qctl4onx"; minPseudo = "24mprbt056"; symp = "fzegt6t3y1"; gob = "wtc1en8xqx" }
)

# The existing data block ends at row 137 (subj_num 314); new rows start at 138.
$startRow = 138

# Values get registered into the workbook's shared-string table in the order
# the source data was generated (alphabetically by test_code: gob, minPseudo,
# minWord, symp) even though the rows themselves are laid out in column order
# (minWord, minPseudo, symp, gob). Fill column H in that alphabetical order
# per subject block first so new strings land in the same order as the source.
$alphaTestCodes = @("gob", "minPseudo", "minWord", "symp")
for ($s = 0; $s -lt $newSubjects.Count; $s++) {
    $subject = $newSubjects[$s]
    $blockRow = $startRow + ($s * 4)
    foreach ($tc in $alphaTestCodes) {
        $offset = [Array]::IndexOf($rowTestCodes, $tc)
        $ws.Cells.Item($blockRow + $offset, 8).Value = $subject[$tc]
    }
}

# Now fill in the subj_num (A) and test_code (G) columns, and confirm/overwrite
# column H, in plain row order for every new row.
$row = $startRow
foreach ($subject in $newSubjects) {
    for ($i = 0; $i -lt 4; $i++) {
        $tc = $rowTestCodes[$i]
        $ws.Cells.Item($row, 1).Value = $subject.Subj   # column A: subj_num
        $ws.Cells.Item($row, 7).Value = $tc             # column G: test_code
        $ws.Cells.Item($row, 8).Value = $subject[$tc]    # column H: value
        $row = $row + 1
    }
}

# Scroll/select to show the newly added bottom of the table, same as the
# author did after pasting in the new rows.
$ws.Range("H164").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 131
$win.ScrollColumn = 1
